$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "5.48", "0.0₃0743") are preserved exactly as text
# instead of being coerced into floating point numbers by Excel.
$ws.Columns.Item(4).NumberFormat = "@"
$ws.Columns.Item(5).NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.615.77"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3
$ws.Range("D3").Value = "3.311.98"
$ws.Range("E3").Value = "  +4.76%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "607.61"
$ws.Range("E5").Value = "  +2.82%  "

# Row 6
$ws.Range("D6").Value = "142.12"
$ws.Range("E6").Value = "  +2.48%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.307.38"
$ws.Range("E8").Value = "  +4.67%  "

# Row 9
$ws.Range("E9").Value = "  +0.29%  "

# Row 10
$ws.Range("E10").Value = "  +2.41%  "

# Row 11
$ws.Range("D11").Value = "5.48"
$ws.Range("E11").Value = "  +3.63%  "

# Row 12
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  +2.05%  "

# Row 13
$ws.Range("E13").Value = "  +0.75%  "

# Row 14
$ws.Range("D14").Value = "34.59"
$ws.Range("E14").Value = "  +0.94%  "

# Row 15
$ws.Range("D15").Value = "3.860.00"
$ws.Range("E15").Value = "  +4.81%  "

# Row 16
$ws.Range("E16").Value = "  +0.23%  "

# Row 17
$ws.Range("D17").Value = "3.316.00"
$ws.Range("E17").Value = "  +5.04%  "

# Row 18
$ws.Range("D18").Value = "63.745.60"
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("E19").Value = "  +2.44%  "

# Row 20
$ws.Range("D20").Value = "479.54"
$ws.Range("E20").Value = "  +0.53%  "

# Row 21
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +0.40%  "

# Row 22
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").Value = "  +4.23%  "

# Row 23
$ws.Range("D23").Value = "8.12"
$ws.Range("E23").Value = "  +5.13%  "

# Row 24
$ws.Range("D24").Value = "13.67"
$ws.Range("E24").Value = "  +4.96%  "

# Row 25
$ws.Range("D25").Value = "84.77"
$ws.Range("E25").Value = "  +0.12%  "

# Row 26
$ws.Range("E26").Value = "  +0.21%  "

# Row 27
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +2.14%  "

# Row 28
$ws.Range("D28").Value = "7.34"
$ws.Range("E28").Value = "  +1.74%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("D30").Value = "8.13"

# Row 31
$ws.Range("D31").Value = "2.15"
$ws.Range("E31").Value = "  +0.79%  "

# Row 32
$ws.Range("D32").Value = "28.81"
$ws.Range("E32").Value = "  +6.68%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  +3.01%  "

# Row 36
$ws.Range("E36").Value = "  +3.37%  "

# Row 37
$ws.Range("D37").Value = "52.67"
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  +5.76%  "

# Row 39
$ws.Range("D39").Value = "0.0400"
$ws.Range("E39").Value = "  +2.48%  "

# Row 40
$ws.Range("D40").Value = "431.42"
$ws.Range("E40").Value = "  +2.29%  "

# Row 41
$ws.Range("D41").Value = "3.071.31"
$ws.Range("E41").Value = "  +4.73%  "

# Row 42
$ws.Range("E42").Value = "  -0.59%  "

# Row 43
$ws.Range("E43").Value = "  +0.37%  "

# Row 44
$ws.Range("D44").Value = "0.116"
$ws.Range("E44").Value = "  +3.40%  "

# Row 45
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("E46").Value = "  +2.64%  "

# Row 47
$ws.Range("D47").Value = "26.30"
$ws.Range("E47").Value = "  +2.94%  "

# Row 48
$ws.Range("D48").Value = "35.96"
$ws.Range("E48").Value = "  +11.99%  "

# Row 49
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.114"
$ws.Range("E50").Value = "  +0.59%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "124.72"
$ws.Range("E51").Value = "  +2.91%  "
